{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Title: \"<AI Pathfinding>\" -> \"<AI Pathfinding Assignment 2>\"\n//    Insert \" Assignment 2\" right after the \"AI Pathfinding\" text and\n//    before the closing \">\".\n// ---------------------------------------------------------------------\nconst titleResults = body.search(\"AI Pathfinding\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\" Assignment 2\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Subtitle date: \"XX.XX.20XX\" -> \"12.05.2021\"\n// ---------------------------------------------------------------------\nconst dateResults = body.search(\"XX.XX.20XX\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"12.05.2021\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) \"Your Student ID\" -> \"880606479\"\n// ---------------------------------------------------------------------\nconst idResults = body.search(\"Your Student ID\", { matchCase: true });\nidResults.load(\"text\");\nawait context.sync();\nif (idResults.items.length > 0) {\n  idResults.items[0].insertText(\"880606479\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) Replace the two \"Inheritance in AI scripts...\" / \"Funtional vs OOP\"\n//    bullet list paragraphs with a single plain paragraph of new text,\n//    and add one extra blank paragraph after it.\n// ---------------------------------------------------------------------\nconst bulletResults = body.search(\"Inheritance in AI scripts\", { matchCase: true });\nbulletResults.load(\"text\");\nawait context.sync();\n\nif (bulletResults.items.length > 0) {\n  const firstBulletPara = bulletResults.items[0].paragraphs.getFirst();\n  firstBulletPara.load(\"text\");\n  await context.sync();\n\n  const secondBulletPara = firstBulletPara.getNext();\n  secondBulletPara.load(\"text\");\n  const anchorPara = secondBulletPara.getNext(); // blank, non-list paragraph right after\n  anchorPara.load(\"text\");\n  await context.sync();\n\n  const newText =\n    \"Object-Oriented programming is used to create ai from Object-Oriented\\u2019s ability to use abstract data types, Class objects and Instance methods. The use of inheritance can also be extremely helpful when developing an AI. You can reference other scripts and use their sections to help develop other scripts in a non-complex manor \";\n\n  // Insert the new plain paragraph and a new blank paragraph before the\n  // anchor (which keeps its own original formatting untouched).\n  anchorPara.insertParagraph(newText, Word.InsertLocation.before);\n  anchorPara.insertParagraph(\"\", Word.InsertLocation.before);\n  await context.sync();\n\n  // Remove the two old list paragraphs.\n  firstBulletPara.delete();\n  secondBulletPara.delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Title: \"<AI Pathfinding>\" -> \"<AI Pathfinding Assignment 2>\"\n#    Insert \" Assignment 2\" right after \"AI Pathfinding\" and before the\n#    closing \">\".\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$null = $rng.Find.Execute(\"AI Pathfinding\")\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertAfter(\" Assignment 2\")\n\n# ---------------------------------------------------------------------\n# 2) Subtitle date: \"XX.XX.20XX\" -> \"12.05.2021\"\n# ---------------------------------------------------------------------\n$rng2 = $d.Content\n$null = $rng2.Find.Execute(\"XX.XX.20XX\")\n$rng2.Text = \"12.05.2021\"\n\n# ---------------------------------------------------------------------\n# 3) \"Your Student ID\" -> \"880606479\"\n# ---------------------------------------------------------------------\n$rng3 = $d.Content\n$null = $rng3.Find.Execute(\"Your Student ID\")\n$rng3.Text = \"880606479\"\n\n# ---------------------------------------------------------------------\n# 4) Replace the two \"Inheritance in AI scripts...\" / \"Funtional vs OOP\"\n#    bullet list paragraphs with a single plain paragraph of new text,\n#    and add one extra blank paragraph after it.\n# ---------------------------------------------------------------------\n$rng4 = $d.Content\n$null = $rng4.Find.Execute(\"Inheritance in AI scripts\")\n$firstBulletPara = $rng4.Paragraphs.Item(1)\n$firstBulletIndex = $firstBulletPara.Range.Paragraphs.Item(1).Index\n\n$paras = $d.Paragraphs\n$firstIdx = $firstBulletPara.Range.Information(3)  # wdActiveEndAdjustedPageNumber placeholder (unused)\n\n# Locate paragraph indices directly via the Paragraphs collection by index.\n$allParas = $d.Paragraphs\n$total = $allParas.Count\n$firstBulletIdxFound = -1\nfor ($i = 1; $i -le $total; $i++) {\n    if ($allParas.Item($i).Range.Text -like \"Inheritance in AI scripts*\") {\n        $firstBulletIdxFound = $i\n        break\n    }\n}\n\n$secondBulletIdxFound = $firstBulletIdxFound + 1\n$anchorIdxFound = $firstBulletIdxFound + 2   # the blank, non-list paragraph right after\n\n$anchorPara = $allParas.Item($anchorIdxFound)\n$anchorRange = $anchorPara.Range\n$anchorRange.InsertParagraphBefore()   # new blank paragraph (placed right before the anchor)\n$anchorRange.InsertParagraphBefore()   # new content paragraph (currently blank)\n\n$newContentPara = $d.Paragraphs.Item($firstBulletIdxFound)\n$newContentPara.Range.Text = \"Object-Oriented programming is used to create ai from Object-Oriented\" + [char]0x2019 + \"s ability to use abstract data types, Class objects and Instance methods. The use of inheritance can also be extremely helpful when developing an AI. You can reference other scripts and use their sections to help develop other scripts in a non-complex manor \"\n\n# Delete the two original list paragraphs (now at firstBulletIdxFound+1 / +2\n# since the new content paragraph pushed them down by one).\n$parasAfterInsert = $d.Paragraphs\n$parasAfterInsert.Item($firstBulletIdxFound + 2).Range.Delete()\n$parasAfterInsert.Item($firstBulletIdxFound + 1).Range.Delete()\n"}
